$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3544.3
$ws.Range("I86").Value = 2593
$ws.Range("K86").Value = 2593
$ws.Range("M86").Value = -1470

$ws.Range("H89").Value = 3544.3
$ws.Range("I89").Value = 2593
$ws.Range("K89").Value = 12965
$ws.Range("M89").Value = -7349

$ws.Range("H92").Value = 568.8333
$ws.Range("I92").Value = 632.6
$ws.Range("K92").Value = 632.6
$ws.Range("M92").Value = 615.4

$ws.Range("H99").Value = 1498.2
$ws.Range("I99").Value = 1432.3334
$ws.Range("K99").Value = 4297.0002
$ws.Range("M99").Value = -2799.0002

$ws.Range("H107").Value = 1855.5
$ws.Range("I107").Value = 835.2857
$ws.Range("K107").Value = 835.2857
$ws.Range("M107").Value = 1084.7143

$ws.Range("H113").Value = 4946.077
$ws.Range("J113").Value = 6914.2856
$ws.Range("L113").Value = 6914.2856
$ws.Range("N113").Value = -13422.2856

$ws.Range("H116").Value = 6438.8335
$ws.Range("I116").Value = 5950
$ws.Range("J116").Value = 6683.25
$ws.Range("K116").Value = 5950
$ws.Range("L116").Value = 6683.25
$ws.Range("M116").Value = -2508
$ws.Range("N116").Value = -13567.25

$ws.Range("H131").Value = 2078
$ws.Range("I131").Value = 2343.5
$ws.Range("J131").Value = 1812.5
$ws.Range("K131").Value = 7030.5
$ws.Range("L131").Value = 5437.5
$ws.Range("M131").Value = -1990.5
$ws.Range("N131").Value = -15517.5

$ws.Range("H132").Value = 125004730
$ws.Range("I132").Value = 142861700
$ws.Range("K132").Value = 428585100
$ws.Range("M132").Value = -428582570

$ws.Range("H137").Value = 2596.077
$ws.Range("J137").Value = 3212.4285
$ws.Range("L137").Value = 9637.2855
$ws.Range("N137").Value = -14737.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1650
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 1728
$ws.Range("I61").Value = 1742
$ws.Range("K61").Value = 1742
$ws.Range("M61").Value = -1530

$ws.Range("H74").Value = 2777.75
$ws.Range("I74").Value = 3003.6667
$ws.Range("J74").Value = 2100
$ws.Range("K74").Value = 3003.6667
$ws.Range("L74").Value = 2100
$ws.Range("M74").Value = -2129.6667
$ws.Range("N74").Value = -3848

$ws.Range("H77").Value = 2777.75
$ws.Range("I77").Value = 3003.6667
$ws.Range("J77").Value = 2100
$ws.Range("K77").Value = 15018.3335
$ws.Range("L77").Value = 10500
$ws.Range("M77").Value = -10650.3335
$ws.Range("N77").Value = -19236

$ws.Range("H122").Value = 3145.111
$ws.Range("I122").Value = 3145.111
$ws.Range("K122").Value = 9435.332999999999
$ws.Range("M122").Value = -6985.332999999999

$ws.Range("H132").Value = 4515.357
$ws.Range("I132").Value = 4969.75
$ws.Range("K132").Value = 14909.25
$ws.Range("M132").Value = -12379.25

$ws.Range("H136").Value = 1728
$ws.Range("I136").Value = 1742
$ws.Range("K136").Value = 5226
$ws.Range("M136").Value = -2676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1034.5714
$ws.Range("I37").Value = 1213
$ws.Range("J37").Value = 796.6667
$ws.Range("K37").Value = 1213
$ws.Range("L37").Value = 796.6667
$ws.Range("M37").Value = -1076
$ws.Range("N37").Value = -1070.6667

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H134").Value = 4807.8667
$ws.Range("I134").Value = 4835.724
$ws.Range("K134").Value = 14507.172
$ws.Range("M134").Value = -11972.172

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 465
$ws.Range("I22").Value = 112.5
$ws.Range("K22").Value = 112.5
$ws.Range("M22").Value = 237.5

$ws.Range("H25").Value = 1461.6
$ws.Range("I25").Value = 1461.6
$ws.Range("K25").Value = 1461.6
$ws.Range("M25").Value = -1287.6

$ws.Range("H74").Value = 39221.668
$ws.Range("J74").Value = 39221.668
$ws.Range("L74").Value = 39221.668
$ws.Range("N74").Value = -40969.668

$ws.Range("H77").Value = 39221.668
$ws.Range("J77").Value = 39221.668
$ws.Range("L77").Value = 117665.004
$ws.Range("N77").Value = -126401.004

$ws.Range("H103").Value = 6883.3335
$ws.Range("I103").Value = 6883.3335
$ws.Range("K103").Value = 6883.3335
$ws.Range("M103").Value = -5711.3335

$ws.Range("H107").Value = 479.7143
$ws.Range("J107").Value = 581.6
$ws.Range("L107").Value = 581.6
$ws.Range("N107").Value = -4421.6

$ws.Range("H132").Value = 3098.625
$ws.Range("J132").Value = 4999.5
$ws.Range("L132").Value = 14998.5
$ws.Range("N132").Value = -20058.5

$ws.Range("H134").Value = 2734
$ws.Range("I134").Value = 1958.5
$ws.Range("K134").Value = 5875.5
$ws.Range("M134").Value = -3340.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45506.52
$ws.Range("J4").Value = 2437.5
$ws.Range("L4").Value = 7312.5
$ws.Range("N4").Value = -7536.5

$ws.Range("H18").Value = 614
$ws.Range("I18").Value = 614
$ws.Range("K18").Value = 1842
$ws.Range("M18").Value = -1673

$ws.Range("H37").Value = 97139.28999999999
$ws.Range("J37").Value = 97139.28999999999
$ws.Range("L37").Value = 291417.87
$ws.Range("N37").Value = -291641.87

$ws.Range("H39").Value = 2199.4
$ws.Range("J39").Value = 3998.5
$ws.Range("L39").Value = 11995.5
$ws.Range("N39").Value = -12583.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4815.385
$ws.Range("I80").Value = 3150
$ws.Range("J80").Value = 5555.5557
$ws.Range("K80").Value = 3150
$ws.Range("L80").Value = 5555.5557
$ws.Range("M80").Value = -2152
$ws.Range("N80").Value = -7551.5557

$ws.Range("H83").Value = 4815.385
$ws.Range("I83").Value = 3150
$ws.Range("J83").Value = 5555.5557
$ws.Range("K83").Value = 15750
$ws.Range("L83").Value = 27777.7785
$ws.Range("M83").Value = -10758
$ws.Range("N83").Value = -37761.7785

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 843.9
$ws.Range("J55").Value = 474
$ws.Range("L55").Value = 474
$ws.Range("N55").Value = -820

$ws.Range("H132").Value = 15567.6
$ws.Range("I132").Value = 17990.385
$ws.Range("K132").Value = 53971.155
$ws.Range("M132").Value = -51441.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 974.4
$ws.Range("I122").Value = 974.4
$ws.Range("K122").Value = 2923.2
$ws.Range("M122").Value = -473.1999999999998

$ws.Range("H136").Value = 2470.611
$ws.Range("I136").Value = 2342
$ws.Range("J136").Value = 3499.5
$ws.Range("K136").Value = 7026
$ws.Range("L136").Value = 10498.5
$ws.Range("M136").Value = -4476
$ws.Range("N136").Value = -15598.5
